# Weekly refresh of fruit/vegetable price data: the rows for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Berenjena" get
# reshuffled (each row's D/J/K/L/M/O/P values move to a different row).
# Row 1 (header), row 13 and row 16 are left untouched.
#
# Note: this runtime's COM .Value getter/setter does not reliably
# round-trip scalars through intermediate PowerShell variables, so we
# consistently use .Value2 instead (works for numbers, dates and
# strings alike).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the columns that move:
#   D = Fecha, J = Volumen, K = Precio minimo, L = Precio maximo,
#   M = Precio promedio ponderado, O = Origen, P = Precio $/Kg
$rows = @(2,3,4,5,6,7,8,9,10,11,12,14,15)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Mapping: destination row -> source row (where the new data for that
# row comes from in the original layout).
$map = @{
    2  = 12
    3  = 5
    4  = 7
    5  = 9
    6  = 8
    7  = 2
    8  = 6
    9  = 11
    10 = 3
    11 = 14
    12 = 15
    14 = 10
    15 = 4
}

foreach ($dest in $map.Keys) {
    $src = $map[$dest]
    $data = $snapshot[$src]
    $ws.Cells.Item($dest, 4).Value2  = $data.D
    $ws.Cells.Item($dest, 10).Value2 = $data.J
    $ws.Cells.Item($dest, 11).Value2 = $data.K
    $ws.Cells.Item($dest, 12).Value2 = $data.L
    $ws.Cells.Item($dest, 13).Value2 = $data.M
    $ws.Cells.Item($dest, 15).Value2 = $data.O
    $ws.Cells.Item($dest, 16).Value2 = $data.P
}
